# Sprint_Backlog.xlsx - "Backlog" sheet gained a new backlog item
# ("Manage group", priority 3) inserted right before the existing
# "Dedicated client" row (old row 38), pushing the following rows down
# by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# Insert a new row at row 38; this shifts the old rows 38-41 down to 39-42
# and keeps everything above row 38 untouched.
$ws.Rows.Item(38).Insert()

$ws.Range("A38").Value = 3
$ws.Range("B38").Value = "Manage group"

# Reflect the new view state (scrolled down a bit further, new selection)
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H40").Select()
